$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 13 (pushes rows 13-25 down to 14-26)
$ws.Rows("13:13").Insert()

# --- Row 13 (new blank row): label cells left/middle aligned, no border ---
$ws.Range("B13:E13").HorizontalAlignment = -4131
$ws.Range("B13:E13").VerticalAlignment = -4108

# --- Row 15 (previously row 14): header row "Ady | Mukdar | Baha | Jemi" ---
# shift Mukdar/Baha/Jemi left one column (drop the "Kody" column) and merge A15:B15 for "Ady"
$mukdar = $ws.Range("B15").Value2
$ws.Range("C15").Value = $mukdar
$ws.Range("B15").Value = ""
$ws.Range("A15:B15").Merge()
$ws.Range("A15:B15").HorizontalAlignment = -4108
$ws.Range("A15:B15").VerticalAlignment = -4108

# --- Row 17 (previously row 16): clear the "180.80 m" value next to "Jemi:" ---
$ws.Range("D17").Value = ""

Write-Host "edit applied"
